$d = $word.ActiveDocument

# Locate the single paragraph that holds the "Chi phi van hanh, quan ly, hanh
# chinh" cost line (under "Uoc luong gia thanh") and rewrite it in place to
# "Chi phi kinh doanh, quang cao, tiep thi: 1.000.000 VND".
#
# Every Find/Execute below is scoped to that paragraph's own Range (re-read
# fresh each time since prior edits shift character offsets) and uses
# wdReplaceOne (1) so only the intended occurrence is touched, never the
# whole story.

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*vận hành, quản lý, hành chính*") {

        # "vận" -> "kinh"
        $r = $p.Range
        $r.Find.Execute("vận", $true, $false, $false, $false, $false, $true, 0, $false, "kinh", 1)

        # first "hành" (now right after "kinh") -> "doanh"
        $r = $p.Range
        $r.Find.Execute("hành", $true, $false, $false, $false, $false, $true, 0, $false, "doanh", 1)

        # "quản" -> "quảng"
        $r = $p.Range
        $r.Find.Execute("quản", $true, $false, $false, $false, $false, $true, 0, $false, "quảng", 1)

        # "lý" -> "cáo"
        $r = $p.Range
        $r.Find.Execute("lý", $true, $false, $false, $false, $false, $true, 0, $false, "cáo", 1)

        # remaining "hành" -> "tiếp"
        $r = $p.Range
        $r.Find.Execute("hành", $true, $false, $false, $false, $false, $true, 0, $false, "tiếp", 1)

        # "chính" -> "thị"
        $r = $p.Range
        $r.Find.Execute("chính", $true, $false, $false, $false, $false, $true, 0, $false, "thị", 1)

        # "2.000.000 VNĐ" -> "1.000.000 VNĐ"
        $r = $p.Range
        $r.Find.Execute("2.000.000 VNĐ", $true, $false, $false, $false, $false, $true, 0, $false, "1.000.000 VNĐ", 1)

        break
    }
}
